$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.765.41'
$ws.Range("E2").Value = '  -1.24%  '

$ws.Range("D3").Value = '1.600.60'
$ws.Range("E3").Value = '  -1.83%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.80%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.505'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.88%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.248'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0617'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.60'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0837'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.30%  '

$ws.Range("D12").Value = '1.832.97'

$ws.Range("D13").Value = '1.614.17'
$ws.Range("E13").Value = '  -0.78%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.05'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.36%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.529'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.23%  '

$ws.Range("D16").Value = '26.762.86'
$ws.Range("E16").Value = '  -1.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.44'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.47%  '

$ws.Range("D18").Value = '0.0₃0726'
$ws.Range("E18").Value = '  -0.71%  '

$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.03%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '209.27'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.23%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.72'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.48%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.78%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.32'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.19%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.84'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.16'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.69%  '

$ws.Range("E26").Value = '  +1.13%  '

$ws.Range("E27").Value = '  +0.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.112'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.88%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.30'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.88%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0500'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.71%  '

$ws.Range("E31").Value = '  -2.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.24'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.667'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +23.30%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.94'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.10%  '

$ws.Range("D35").Value = '1.312.17'
$ws.Range("E35").Value = '  +0.86%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.51'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.98%  '

$ws.Range("E37").Value = '  -0.40%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0172'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.67%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.818'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.95%  '

$ws.Range("E40").Value = '  -0.01%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.788'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.28%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.18'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.51%  '

$ws.Range("E43").Value = '  -0.20%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.69'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.73%  '

$ws.Range("D45").Value = '1.744.25'
$ws.Range("E45").Value = '  -1.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.95'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.81%  '

$ws.Range("E47").Value = '  +0.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.810'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.11%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0510'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.55%  '

$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0975'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.79%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₇0975'
$ws.Range("E51").Value = '  -8.94%  '
